$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates — force text to preserve the non-numeric /
# pseudo-numeric string formatting used in the source data (dates-like
# "23.188.07", plain decimals "1.001", tiny values "0.00001270", etc.)
$dUpdates = @{
    'D2' = '23.188.07'
    'D3' = '1.609.14'
    'D4' = '1.001'
    'D5' = '0.9997'
    'D6' = '302.89'
    'D7' = '0.3775'
    'D8' = '0.3664'
    'D9' = '49.14'
    'D10' = '1.001'
    'D11' = '1.274'
    'D12' = '0.08093'
    'D13' = '23.12'
    'D14' = '6.613'
    'D15' = '7.587'
    'D16' = '0.00001270'
    'D17' = '1.605.77'
    'D18' = '91.74'
    'D19' = '0.06792'
    'D21' = '6.597'
    'D22' = '0.9993'
    'D23' = '13.11'
    'D24' = '23.214.44'
    'D25' = '2.361'
    'D26' = '2.926'
    'D27' = '21.14'
    'D28' = '150.69'
    'D29' = '5.252'
    'D30' = '132.59'
    'D31' = '2.415'
    'D32' = '6.976'
    'D33' = '1.787.21'
    'D34' = '0.9795'
    'D35' = '0.07757'
    'D36' = '0.02787'
    'D37' = '6.308'
    'D38' = '0.2556'
    'D39' = '10.11'
    'D40' = '0.08872'
    'D41' = '1.401'
    'D42' = '0.7172'
    'D43' = '12.79'
    'D44' = '15.98'
    'D45' = '0.6628'
    'D46' = '2.307'
    'D47' = '0.9985'
    'D48' = '3.984'
    'D49' = '0.08031'
    'D50' = '131.68'
    'D51' = '1.171'
}

foreach ($addr in $dUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$addr]
    $cell.Style = "Normal"
}

# Column E (Volume(1h)) updates — these are already non-numeric strings
# (padded with spaces and a trailing "%"), so a direct Value assignment
# keeps them as text.
$eUpdates = @{
    'E3' = '  -2.72%  '
    'E4' = '  -0.05%  '
    'E5' = '  -0.11%  '
    'E6' = '  -2.21%  '
    'E7' = '  -3.25%  '
    'E8' = '  -4.50%  '
    'E9' = '  -4.33%  '
    'E10' = '  +0.00%  '
    'E11' = '  -5.92%  '
    'E12' = '  -4.26%  '
    'E13' = '  -3.45%  '
    'E14' = '  -7.16%  '
    'E15' = '  -3.89%  '
    'E16' = '  -3.55%  '
    'E17' = '  -2.89%  '
    'E18' = '  -3.05%  '
    'E19' = '  -2.89%  '
    'E20' = '  -6.89%  '
    'E21' = '  -4.91%  '
    'E22' = '  -0.17%  '
    'E23' = '  -4.28%  '
    'E24' = '  -3.13%  '
    'E25' = '  -4.92%  '
    'E26' = '  -2.45%  '
    'E27' = '  -4.49%  '
    'E28' = '  -0.73%  '
    'E29' = '  -3.65%  '
    'E30' = '  -4.90%  '
    'E31' = '  -3.02%  '
    'E32' = '  -11.22%  '
    'E33' = '  -2.66%  '
    'E34' = '  -6.04%  '
    'E35' = '  -4.32%  '
    'E36' = '  -6.09%  '
    'E37' = '  -6.95%  '
    'E38' = '  -4.78%  '
    'E39' = '  -7.17%  '
    'E40' = '  -3.11%  '
    'E41' = '  -1.90%  '
    'E42' = '  -5.09%  '
    'E43' = '  -5.07%  '
    'E44' = '  -2.62%  '
    'E45' = '  -4.65%  '
    'E46' = '  -6.09%  '
    'E47' = '  -0.17%  '
    'E48' = '  -2.63%  '
    'E49' = '  -3.21%  '
    'E50' = '  -2.25%  '
    'E51' = '  -3.85%  '
}

foreach ($addr in $eUpdates.Keys) {
    $ws.Range($addr).Value = $eUpdates[$addr]
}
